$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country-code mapping cells in column A with new ISO3 codes / country
# names pulled from the refreshed FAO/FABLE mapping database.
# Row 157: R_OEU -> GRC (Greece)
$ws.Range("A157").Value = "GRC"
# Row 67: R_ASIPAC -> NPL (Nepal)
$ws.Range("A67").Value = "NPL"
# Row 146: R_NEU -> Turkey
$ws.Range("A146").Value = "Turkey"
# Row 154: R_OEU -> Denmark
$ws.Range("A154").Value = "Denmark"

# Update the active selection left over from the author's last save.
[void]$ws.Range("B9").Select()
